# Update cryptocurrency price (D) and volume-change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.338.83'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.606.33'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '510.54'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.63'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.81%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('E8').Value = '  -2.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.618.11'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.70'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.26%  '
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('E13').Value = '  +1.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.063.28'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.322.83'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.59'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.611.09'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.65%  '
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '350.55'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.60'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.14'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.97%  '
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('E25').Value = '  -0.47%  '
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0840'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.35'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -2.60%  '
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '151.10'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.93%  '
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('E34').Value = '  +0.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.99'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('E36').Value = '  -2.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.884'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +5.76%  '
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.844'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.31'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.29%  '
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '295.17'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -6.10%  '
$ws.Range('E43').Value = '  -3.38%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('E45').Value = '  +0.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0555'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.85'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.88'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('E49').Value = '  -1.12%  '
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.000.69'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.69%  '
